$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking price text (e.g. trailing zeros like "1.00")
# must be preserved literally as text instead of being auto-converted to a number.
$textCells = @("D5", "D6", "D8", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D31", "D32", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.468.89"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.518.81"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "539.15"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "139.89"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").Value = "2.523.92"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "2.967.52"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "23.50"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "59.370.95"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.517.14"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "326.56"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "63.29"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").Value = "0.425"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "7.83"
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").Value = "6.95"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("D30").Value = "0.0₃0783"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "165.39"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -6.65%  "
$ws.Range("D36").Value = "18.55"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "36.97"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "3.70"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "0.813"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "5.27"
$ws.Range("E42").Value = "  -6.37%  "
$ws.Range("D43").Value = "280.68"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "10.85"
$ws.Range("D47").Value = "0.0934"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "123.86"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "0.0226"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "17.90"
$ws.Range("E51").Value = "  -1.99%  "

# Restore the default (Normal) cell style now that the text values are set,
# so only the displayed text changes and no new number formatting lingers.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
